$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 502 (shifts old rows 502-520 down to 503-521),
# inheriting formatting from the surrounding rows (matches the target diff,
# which adds a "Warring-States Japan Battle Data" entry between the existing
# "Water-Related Intrastate Conflict and Cooperation" row and the
# "Varieties of Party Identity and Organization" row).
$ws.Rows("502:502").Insert()

# Fill in the new row's cell values. The order below matches the order the
# new shared strings appear in the target workbook (name, doi, file_zip,
# link, topics) so the resulting shared-string table lines up with the diff.
$ws.Range("A502").Value = "Warring-States Japan Battle Data"
$ws.Range("Z502").Value = "10.1080/03050629.2023.2149514"
$ws.Range("V502").Value = "https://blogs.gwu.edu/nick_anderson/files/2022/12/WSJBData.zip"
$ws.Range("C502").Value = "https://blogs.gwu.edu/nick_anderson/data/"
$ws.Range("D502").Value = "battles"

$ws.Range("B502").Value = "international relations"
$ws.Range("E502").Value = "JP"
$ws.Range("F502").Value = 0
$ws.Range("G502").Value = 1
$ws.Range("H502").Value = 0
$ws.Range("I502").Value = 0
$ws.Range("J502").Value = 0
$ws.Range("K502").Value = 1467
$ws.Range("L502").Value = 1600
$ws.Range("M502").Value = "online"
$ws.Range("N502").Value = "no"
$ws.Range("O502").Value = 1
$ws.Range("X502").Value = "year"
$ws.Range("AB502").Value = 20221215

# Add the two hyperlinks that point at the new row's URLs (the zip download
# in V502 and the landing page link in C502).
$ws.Hyperlinks.Add($ws.Range("V502"), "https://blogs.gwu.edu/nick_anderson/files/2022/12/WSJBData.zip")
$ws.Hyperlinks.Add($ws.Range("C502"), "https://blogs.gwu.edu/nick_anderson/data/")

# Restore the active selection to match the new row layout.
$ws.Range("A503").Select()

Write-Output "done"
